# The "draft" status indicator was previously an inline picture (an SVG
# badge pulled from rfc.unprotocols.org). Replace it with the plain text
# "draft" so the sentence reads as ordinary text instead of embedding an
# image.
$d = $word.ActiveDocument

$shp = $d.InlineShapes.Item(1)
$shpRange = $shp.Range
$insertAt = $shpRange.Start
$shp.Delete()

$target = $d.Range($insertAt, $insertAt)
$target.InsertAfter("draft")

Write-Output "Replaced draft status image with literal text."
